# Weekly driver report update for 2025-04-21
# Updates Critical Minutes (C) and Good Roaming Calculation (D) figures in the
# "Bad Drivers" section, plus Total Samples (B) figures in the
# "Good Drivers" section, on the active worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Bad Drivers table (rows 3-12) ---

# Row 3: Intel(R) Wi-Fi 6 AX201 160MHz - 22.120.0.3
$ws.Range("C3").Value = 2431
$ws.Range("D3").Value = 86.3

# Row 4: Intel(R) Dual Band Wireless-AC 8265 - 20.70.24.1
$ws.Range("C4").Value = 671

# Row 5: Intel(R) Wi-Fi 6 AX201 160MHz - 22.40.0.7
$ws.Range("C5").Value = 507
$ws.Range("D5").Value = 93.90000000000001

# Row 6: Intel(R) Wi-Fi 6 AX201 160MHz - 23.20.1.1
$ws.Range("C6").Value = 626
$ws.Range("D6").Value = 93.90000000000001

# Row 7: Intel(R) Wi-Fi 6 AX201 160MHz - 22.100.0.3
$ws.Range("C7").Value = 280
$ws.Range("D7").Value = 96.09999999999999

# Row 8: Intel(R) Dual Band Wireless-AC 7265 - 19.51.12.3
$ws.Range("C8").Value = 1105
$ws.Range("D8").Value = 96.90000000000001

# Row 9: Intel(R) Dual Band Wireless-AC 8265 - 20.70.17.1
$ws.Range("C9").Value = 189

# Row 10: Intel(R) Dual Band Wireless-AC 8265 - 20.70.3.3
$ws.Range("C10").Value = 111
$ws.Range("D10").Value = 97.5

# Row 11: Intel(R) Dual Band Wireless-AC 8265 - 20.50.3.3
$ws.Range("C11").Value = 692
$ws.Range("D11").Value = 98.2

# Row 12: Totals
$ws.Range("C12").Value = 6612

# --- Good Drivers table (rows 20-39), Total Samples column B ---

# Row 20: Intel(R) Dual Band Wireless-AC 7265 - 19.50.1.6
$ws.Range("B20").Value = 29731

# Row 22: Intel(R) Wi-Fi 6 AX201 160MHz - 23.100.0.4
$ws.Range("B22").Value = 449371

# Row 24: Intel(R) Wi-Fi 6 AX201 160MHz - 22.80.0.9
$ws.Range("B24").Value = 77999

# Row 29: Intel(R) Dual Band Wireless-AC 7265 - 19.51.30.1
$ws.Range("B29").Value = 205276

# Row 30: Intel(R) Dual Band Wireless-AC 7265 - 19.51.29.1
$ws.Range("B30").Value = 40211

# Row 34: Intel(R) Dual Band Wireless-AC 8265 - 20.70.12.5
$ws.Range("B34").Value = 144782

# Row 39: Intel(R) Dual Band Wireless-AC 7265 - 19.51.14.1
$ws.Range("B39").Value = 122297
